$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Calibrated Results")

$ws.Cells.Item(2, 8).Value = 3.319427019503336
$ws.Cells.Item(3, 8).Value = 0.6151083434800566
$ws.Cells.Item(4, 8).Value = 8.190568674965512
$ws.Cells.Item(5, 8).Value = 3.506075870118107
$ws.Cells.Item(6, 8).Value = 13.40894031935255
$ws.Cells.Item(7, 8).Value = 7.815171439159344
$ws.Cells.Item(8, 8).Value = 4.046835525496548
$ws.Cells.Item(9, 8).Value = 10.86694941796445
$ws.Cells.Item(10, 8).Value = 2.25353956114437
$ws.Cells.Item(11, 8).Value = 4.57464507516263
$ws.Cells.Item(12, 8).Value = 4.958440449276153
$ws.Cells.Item(13, 8).Value = 2.966833537208002
$ws.Cells.Item(14, 8).Value = 2.071089779683275
$ws.Cells.Item(15, 8).Value = 2.245141422837015
$ws.Cells.Item(16, 8).Value = 1.72059583196478
$ws.Cells.Item(17, 8).Value = 5.296919401642119
$ws.Cells.Item(18, 8).Value = 3.136977238042455
$ws.Cells.Item(19, 8).Value = 7.28642677913343
$ws.Cells.Item(20, 8).Value = 7.468294348486022
$ws.Cells.Item(21, 8).Value = 7.4724934176397
$ws.Cells.Item(22, 8).Value = 7.46257795686401
$ws.Cells.Item(23, 8).Value = 4.76397567246253
$ws.Cells.Item(24, 8).Value = 4.228350196597812
$ws.Cells.Item(25, 5).Value = 1
$ws.Cells.Item(25, 6).Value = 0
$ws.Cells.Item(25, 8).Value = 13.0605459062109
$ws.Cells.Item(26, 8).Value = 9.443687196047755
$ws.Cells.Item(27, 8).Value = 2.594118048086961
$ws.Cells.Item(28, 8).Value = 5.491966390564245
$ws.Cells.Item(29, 8).Value = 6.396108286396327
$ws.Cells.Item(30, 8).Value = 4.952141845545636
$ws.Cells.Item(31, 8).Value = 3.701705071148523
$ws.Cells.Item(32, 8).Value = 5.487185109302064
$ws.Cells.Item(33, 8).Value = 7.810972370005667
$ws.Cells.Item(34, 6).Value = 0
$ws.Cells.Item(34, 8).Value = 13.59558916996732
$ws.Cells.Item(35, 8).Value = 2.425491669721272
$ws.Cells.Item(36, 8).Value = 5.485085574725225
$ws.Cells.Item(37, 8).Value = 8.006019358927578
$ws.Cells.Item(38, 8).Value = 7.278028640826077
$ws.Cells.Item(39, 8).Value = 3.883572640501329
$ws.Cells.Item(42, 8).Value = 8.006019358927578
$ws.Cells.Item(43, 8).Value = 3.506075870118107
$ws.Cells.Item(44, 8).Value = 6.023392797275285
$ws.Cells.Item(45, 8).Value = 6.557853848923425
$ws.Cells.Item(46, 8).Value = 6.93173376226126
$ws.Cells.Item(47, 8).Value = 6.01861151601332
$ws.Cells.Item(48, 8).Value = 6.205260366627877
$ws.Cells.Item(49, 8).Value = 7.473075629747989
$ws.Cells.Item(50, 8).Value = 6.935350619306432
$ws.Cells.Item(51, 8).Value = 9.451503122246818
$ws.Cells.Item(52, 8).Value = 3.346429842168062
$ws.Cells.Item(53, 8).Value = 4.041701345983039
$ws.Cells.Item(54, 8).Value = 5.491384178455742
$ws.Cells.Item(55, 8).Value = 7.4724934176397
$ws.Cells.Item(56, 8).Value = 3.158263669084955
$ws.Cells.Item(57, 8).Value = 2.784383755747121
$ws.Cells.Item(58, 8).Value = 13.40894031935255
$ws.Cells.Item(60, 8).Value = 1.53184744677317
$ws.Cells.Item(61, 8).Value = 9.605432758574853
$ws.Cells.Item(62, 8).Value = 4.424914507988059
$ws.Cells.Item(63, 8).Value = 6.020711050590158
$ws.Cells.Item(65, 8).Value = 5.303218005372635
$ws.Cells.Item(66, 8).Value = 5.487185109302064
$ws.Cells.Item(67, 8).Value = 2.252957349035867
$ws.Cells.Item(69, 8).Value = 1.368002349669447
$ws.Cells.Item(70, 8).Value = 6.024910119743836
$ws.Cells.Item(71, 8).Value = 1.711615481548923
$ws.Cells.Item(72, 8).Value = 2.975231675515356
$ws.Cells.Item(73, 8).Value = 5.491966390564245
$ws.Cells.Item(74, 8).Value = 5.674416172025126
$ws.Cells.Item(75, 8).Value = 7.817270973736181
$ws.Cells.Item(76, 8).Value = 8.727711473298779
$ws.Cells.Item(77, 8).Value = 3.693306932841383
$ws.Cells.Item(78, 8).Value = 7.278028640826077
$ws.Cells.Item(79, 8).Value = 6.053430264876683
$ws.Cells.Item(80, 8).Value = 7.817270973736181
$ws.Cells.Item(82, 8).Value = 2.056393037645405
$ws.Cells.Item(83, 8).Value = 6.213658504935232
$ws.Cells.Item(84, 8).Value = 4.419780328474551
$ws.Cells.Item(85, 8).Value = 6.026427442212171
$ws.Cells.Item(86, 8).Value = 8.723512404145101
$ws.Cells.Item(87, 8).Value = 7.284327244556591
$ws.Cells.Item(88, 8).Value = 6.207359901204716
$ws.Cells.Item(89, 8).Value = 1.34099952700472
$ws.Cells.Item(90, 8).Value = 13.40894031935255
$ws.Cells.Item(91, 8).Value = 5.485667786833729
$ws.Cells.Item(92, 8).Value = 6.931151550152755
$ws.Cells.Item(93, 8).Value = 6.906830474173371
$ws.Cells.Item(94, 8).Value = 3.161880526130128
$ws.Cells.Item(95, 8).Value = 6.747184446223326
$ws.Cells.Item(96, 8).Value = 8.186369605811835
$ws.Cells.Item(97, 8).Value = 7.282227709979753
$ws.Cells.Item(98, 8).Value = 5.675933494493676
$ws.Cells.Item(99, 8).Value = 7.466194813909183
$ws.Cells.Item(100, 8).Value = 6.559953383500049
$ws.Cells.Item(101, 8).Value = 5.145671511999214
$ws.Cells.Item(102, 8).Value = 6.363389072109587
$ws.Cells.Item(104, 8).Value = 6.211558970358393
$ws.Cells.Item(105, 8).Value = 5.676515706601965
$ws.Cells.Item(106, 8).Value = 2.784383755747121
$ws.Cells.Item(107, 8).Value = 4.043800880559878
$ws.Cells.Item(108, 8).Value = 3.507240294334901
$ws.Cells.Item(109, 8).Value = 3.510274939271785
$ws.Cells.Item(110, 8).Value = 4.041701345983039
$ws.Cells.Item(111, 8).Value = 6.213658504935232
$ws.Cells.Item(112, 8).Value = 5.675933494493676
$ws.Cells.Item(113, 8).Value = 2.968933071784841
$ws.Cells.Item(114, 8).Value = 5.142989765314086
$ws.Cells.Item(115, 8).Value = 7.815171439159344
$ws.Cells.Item(116, 8).Value = 1.873361044076021
$ws.Cells.Item(117, 8).Value = 4.952141845545636
$ws.Cells.Item(118, 8).Value = 5.865264091793576
$ws.Cells.Item(119, 8).Value = 5.30951660910315
$ws.Cells.Item(120, 8).Value = 7.466777026017688
$ws.Cells.Item(121, 5).Value = 1
$ws.Cells.Item(121, 6).Value = 0
$ws.Cells.Item(121, 8).Value = 13.0605459062109
$ws.Cells.Item(122, 8).Value = 6.020711050590158
$ws.Cells.Item(123, 8).Value = 3.512374473848623
$ws.Cells.Item(124, 8).Value = 2.782284221170283
$ws.Cells.Item(126, 8).Value = 4.925721234989414
$ws.Cells.Item(127, 8).Value = 0.9878238326010974
$ws.Cells.Item(128, 8).Value = 4.605846966980819
$ws.Cells.Item(129, 8).Value = 4.575227287270919
$ws.Cells.Item(130, 8).Value = 4.228350196597812
$ws.Cells.Item(131, 8).Value = 6.02700965432046
$ws.Cells.Item(132, 8).Value = 2.247240957413855
$ws.Cells.Item(133, 8).Value = 4.950042310968798
$ws.Cells.Item(134, 8).Value = 0.8011749819863256
